# Rename the three logo pictures (the two Pearson-logo footers and the
# BTec-logo first-page header) the same way Word itself renames an
# InlineShape: convert it to a (floating) Shape so .Name becomes
# writable, set the new name, then flip WrapFormat back to inline so it
# re-serialises as <wp:inline> again (matching the original markup
# shape).
#
#   footer (default)    Pearson logo  image2.png -> image1.png
#   footer (first page)  Pearson logo  image2.png -> image1.png
#   header (first page)  BTec logo     image1.jpg -> image2.jpg

$d = $word.ActiveDocument

function Rename-InlineShapeName($inlineShape, $newName) {
    $shp = $inlineShape.ConvertToShape()
    $shp.Name = $newName
    $shp.WrapFormat.Type = 7   # wdWrapInline - restore inline placement
}

$sec = $d.Sections.Item(1)

# First-page header -> BTec logo: image1.jpg -> image2.jpg
$hdrFirst = $sec.Headers.Item(2)   # wdHeaderFooterFirstPage
foreach ($shp in $hdrFirst.Range.InlineShapes) {
    if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
        Rename-InlineShapeName $shp "image2.jpg"
    }
}

# Default footer -> Pearson logo: image2.png -> image1.png
$ftrDefault = $sec.Footers.Item(1)   # wdHeaderFooterPrimary
foreach ($shp in $ftrDefault.Range.InlineShapes) {
    if ($shp.AlternativeText -like "*PearsonLogo.png") {
        Rename-InlineShapeName $shp "image1.png"
    }
}

# First-page footer -> Pearson logo: image2.png -> image1.png
$ftrFirst = $sec.Footers.Item(2)   # wdHeaderFooterFirstPage
foreach ($shp in $ftrFirst.Range.InlineShapes) {
    if ($shp.AlternativeText -like "*PearsonLogo.png") {
        Rename-InlineShapeName $shp "image1.png"
    }
}
